# "Generate Report for Handback"
#
# The handback process completed for the edc538f4-... file: the zh-cn and
# de-de rows now record the handback target/result files and the handback
# timestamp, the overall status flips from "Ready for handoff" to
# "Handed back: in sync with en-US", and the columns that now hold longer
# filenames/status text are widened to fit.

$wb  = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn     = $wb.Worksheets.Item(2)
$wsDeDe     = $wb.Worksheets.Item(3)

$mdFileName   = "edc538f4-45ab-4d48-a06d-db2a608d1a89.md"
$mdHyperlink  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e26b062c424e11f4aecaf381cfb022c59de2e630/e2e/edc538f4-45ab-4d48-a06d-db2a608d1a89.md"
$zhCnXlf      = "edc538f4-45ab-4d48-a06d-db2a608d1a89.b8a770874777c76c4b5256ee553453b2fe30b1f1.zh-cn.xlf"
$deDeXlf      = "edc538f4-45ab-4d48-a06d-db2a608d1a89.b8a770874777c76c4b5256ee553453b2fe30b1f1.de-de.xlf"
$newStatus    = "Handed back: in sync with en-US"
$zhCnHandbackTime = "2016-10-10 07:09:06"
$deDeHandbackTime = "2016-10-10 07:09:24"

# --- Status column: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# (shows up in the Overview sheet for both zh-cn / de-de, and as the Status
# column on each language sheet)
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# --- zh-cn sheet: record the handback target + handback file, link the target file ---
$wsZhCn.Range("I2").Value = $mdFileName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdHyperlink, $null, $null, $mdFileName) | Out-Null
$wsZhCn.Range("J2").Value = $zhCnXlf
$wsZhCn.Range("K2").Value = $zhCnHandbackTime

# --- de-de sheet: record the handback target + handback file, link the target file ---
$wsDeDe.Range("I2").Value = $mdFileName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdHyperlink, $null, $null, $mdFileName) | Out-Null
$wsDeDe.Range("J2").Value = $deDeXlf
$wsDeDe.Range("K2").Value = $deDeHandbackTime

# --- Column widths: widen the columns that now hold longer text ---
# (ColumnWidth snaps to the nearest 1/6 character-width increment in this
# runtime, so these inputs are chosen to land on/near the generated report's
# widths of 29.9777050018311 and 40.)
$wsOverview.Columns.Item(5).ColumnWidth  = 29.16667
$wsOverview.Columns.Item(6).ColumnWidth  = 29.16667

$wsZhCn.Columns.Item(3).ColumnWidth  = 29.16667
$wsZhCn.Columns.Item(9).ColumnWidth  = 39.16667
$wsZhCn.Columns.Item(10).ColumnWidth = 39.16667

$wsDeDe.Columns.Item(3).ColumnWidth  = 29.16667
$wsDeDe.Columns.Item(9).ColumnWidth  = 39.16667
$wsDeDe.Columns.Item(10).ColumnWidth = 39.16667
